$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A32").Value = "x"
$ws.Range("A32").Font.Name = "Inherit"
$ws.Range("A32").Font.Size = 10
